$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking-style "Price" strings (column D) can look like plain numbers
# (e.g. "505.79") or like dotted-thousands numbers (e.g. "56.212.37").
# Excel auto-converts a Range.Value assignment that parses cleanly as a
# number into a numeric cell, which would silently drop formatting like
# the trailing zero in "1.00". Force those through as literal text by
# flipping NumberFormat to "@" for the write, then restoring the original
# style so no stray style index gets left behind.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "56.212.37"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "2.366.92"
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue "D5" "505.79"
$ws.Range("E5").Value = "  +0.54%  "
Set-TextValue "D6" "130.18"
$ws.Range("E6").Value = "  -0.97%  "
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -1.48%  "
$ws.Range("D9").Value = "2.372.89"
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("E11").Value = "  -0.06%  "
Set-TextValue "D12" "4.87"
$ws.Range("E12").Value = "  +7.44%  "
Set-TextValue "D13" "0.324"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("D14").Value = "2.786.32"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").Value = "56.199.71"
$ws.Range("E15").Value = "  -0.82%  "
Set-TextValue "D16" "21.76"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "2.347.69"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("E19").Value = "  -1.50%  "
Set-TextValue "D20" "309.67"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("E22").Value = "  +0.19%  "
Set-TextValue "D23" "1.00"
$ws.Range("E23").Value = "  +0.02%  "
Set-TextValue "D24" "65.55"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("E26").Value = "  -0.49%  "
Set-TextValue "D27" "0.148"
$ws.Range("E27").Value = "  -1.83%  "
Set-TextValue "D28" "7.19"
$ws.Range("E28").Value = "  -2.82%  "
Set-TextValue "D29" "173.09"
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("E31").Value = "  -0.91%  "
Set-TextValue "D32" "5.85"
$ws.Range("E32").Value = "  -0.46%  "
Set-TextValue "D34" "0.996"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  -3.50%  "
Set-TextValue "D36" "17.54"
$ws.Range("E36").Value = "  -1.65%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  -2.94%  "
Set-TextValue "D39" "0.826"
$ws.Range("E39").Value = "  +2.03%  "
Set-TextValue "D40" "36.27"
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("E41").Value = "  -3.17%  "
$ws.Range("E42").Value = "  +0.70%  "
Set-TextValue "D43" "125.41"
$ws.Range("E43").Value = "  -5.28%  "
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("E46").Value = "  -0.70%  "
Set-TextValue "D47" "236.79"
$ws.Range("E47").Value = "  -4.84%  "
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("E49").Value = "  -1.37%  "
Set-TextValue "D50" "16.91"
$ws.Range("E50").Value = "  -0.64%  "
Set-TextValue "D51" "0.952"
$ws.Range("E51").Value = "  +0.13%  "
